$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f3f5a60>),
                (''model'',
                 RandomForestClassifier(max_depth=1, min_samples_leaf=6,
                                        n_estimators=50, random_state=42))])'
$ws.Range('B2').Value = 0.7134215784215785
$ws.Range('C2').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9d10fa90>, ''scaler'': StandardScaler(), ''model__n_estimators'': 50, ''model__min_samples_split'': 2, ''model__min_samples_leaf'': 6, ''model__max_features'': ''sqrt'', ''model__max_depth'': 1, ''model__class_weight'': None}'
$ws.Range('D2').Value = 0.8401853282782578
$ws.Range('E2').Value = 0.5821970640470641
$ws.Range('F2').Value = 0.8108108108108109
$ws.Range('G2').Value = 0.8325774940421917
$ws.Range('H2').Value = 0.5800686507936508
$ws.Range('I2').Value = 0.7142857142857143
$ws.Range('J2').Value = 0.8570638297872339
$ws.Range('K2').Value = 0.6143333333333333
$ws.Range('L2').Value = 0.9375
$ws.Range('M2').Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range('N2').Value = '[1 1 1 1 1 1 0 1 1 1 1 1 1 1 0 0 1 1 1 1 1 1 1 1]'

# Row 3
$ws.Range('A3').Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f5cb9a0>),
                (''model'',
                 RandomForestClassifier(max_depth=1, min_samples_leaf=6,
                                        n_estimators=50, random_state=42))])'
$ws.Range('B3').Value = 0.7107692307692308
$ws.Range('C3').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9d1257c0>, ''scaler'': StandardScaler(), ''model__n_estimators'': 50, ''model__min_samples_split'': 2, ''model__min_samples_leaf'': 6, ''model__max_features'': ''sqrt'', ''model__max_depth'': 1, ''model__class_weight'': None}'
$ws.Range('D3').Value = 0.8329874267787837
$ws.Range('E3').Value = 0.5712759851259851
$ws.Range('F3').Value = 0.8205128205128205
$ws.Range('G3').Value = 0.8304463240857151
$ws.Range('H3').Value = 0.6431464285714286
$ws.Range('I3').Value = 0.6956521739130435
$ws.Range('J3').Value = 0.845191489361702
$ws.Range('K3').Value = 0.5478333333333334
$ws.Range('L3').Value = 1
$ws.Range('M3').Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range('N3').Value = '[1 1 1 1 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 1 1 1 1]'

# Row 4
$ws.Range('A4').Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f6335b0>),
                (''model'',
                 RandomForestClassifier(max_depth=1, max_features=''log2'',
                                        min_samples_leaf=5, min_samples_split=4,
                                        n_estimators=50, random_state=42))])'
$ws.Range('B4').Value = 0.6796503496503496
$ws.Range('C4').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9d13b490>, ''scaler'': None, ''model__n_estimators'': 50, ''model__min_samples_split'': 4, ''model__min_samples_leaf'': 5, ''model__max_features'': ''log2'', ''model__max_depth'': 1, ''model__class_weight'': None}'
$ws.Range('D4').Value = 0.8475767315755198
$ws.Range('E4').Value = 0.577533893883894
$ws.Range('F4').Value = 0.7058823529411765
$ws.Range('G4').Value = 0.8399166819802077
$ws.Range('H4').Value = 0.5925865079365079
$ws.Range('I4').Value = 0.8
$ws.Range('J4').Value = 0.8608222222222223
$ws.Range('K4').Value = 0.5979999999999999
$ws.Range('L4').Value = 0.631578947368421
$ws.Range('M4').Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range('N4').Value = '[0 1 1 1 0 1 0 1 1 0 1 0 0 1 1 1 0 1 1 0 1 0 1 1]'

# Row 5
$ws.Range('A5').Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9d10f0a0>),
                (''model'',
                 RandomForestClassifier(max_depth=1, max_features=''log2'',
                                        min_samples_leaf=4, min_samples_split=4,
                                        n_estimators=50, random_state=42))])'
$ws.Range('B5').Value = 0.7459890109890109
$ws.Range('C5').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9ce09190>, ''scaler'': RobustScaler(), ''model__n_estimators'': 50, ''model__min_samples_split'': 4, ''model__min_samples_leaf'': 4, ''model__max_features'': ''log2'', ''model__max_depth'': 1, ''model__class_weight'': None}'
$ws.Range('D5').Value = 0.8438098805910011
$ws.Range('E5').Value = 0.5974429958929959
$ws.Range('F5').Value = 0.7567567567567568
$ws.Range('G5').Value = 0.8242784681730062
$ws.Range('H5').Value = 0.5800861111111111
$ws.Range('I5').Value = 0.6086956521739131
$ws.Range('J5').Value = 0.8753061224489795
$ws.Range('K5').Value = 0.6426666666666666
$ws.Range('L5').Value = 1
$ws.Range('M5').Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range('N5').Value = '[1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 1]'

# Row 6
$ws.Range('A6').Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9d125e20>),
                (''model'',
                 RandomForestClassifier(max_depth=2, max_features=''log2'',
                                        min_samples_leaf=6, min_samples_split=4,
                                        n_estimators=5, random_state=42))])'
$ws.Range('B6').Value = 0.7533516483516484
$ws.Range('C6').Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9d13b820>, ''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__min_samples_split'': 4, ''model__min_samples_leaf'': 6, ''model__max_features'': ''log2'', ''model__max_depth'': 2, ''model__class_weight'': None}'
$ws.Range('D6').Value = 0.8431748806668244
$ws.Range('E6').Value = 0.6281535797535798
$ws.Range('F6').Value = 0.6470588235294118
$ws.Range('G6').Value = 0.8387768999159489
$ws.Range('H6').Value = 0.6050103174603174
$ws.Range('I6').Value = 0.4782608695652174
$ws.Range('J6').Value = 0.861
$ws.Range('K6').Value = 0.6828333333333333
$ws.Range('L6').Value = 1
$ws.Range('M6').Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range('N6').Value = '[1 1 1 1 1 1 1 1 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1]'
